# The DDAS upload template now accepts a second project number. Insert a
# new "Project Number 2" column right after the existing "Project Number"
# column (i.e. before the old "Sponsor Protocol Number" column), shifting
# every column from the old C onward one position to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting at C pushes the old C:AI block to D:AJ (values/styles carried
# along) and the new column C inherits column B's formatting.
$ws.Range("C1").EntireColumn.Insert()

# Header text for the brand new column.
$ws.Range("C1").Value = "Project Number 2"

# Rows 3 & 4 only ever had a single formatted-but-empty placeholder cell
# in column B ("Sponsor Protocol Number"/blank row filler). The insert
# shifts that placeholder to column C but also leaves a stray formatted
# cell behind in column B - clear it so column B goes back to blank there.
$ws.Range("B3").Clear()
$ws.Range("B4").Clear()

# Column widths: the new "Project Number 2" column (B stays as-is -
# "Project Number"; C is the new column) gets a fixed width, and the
# shifted "Sponsor Protocol Number" column keeps (approximately) its
# previous best-fit width.
$ws.Range("B1").EntireColumn.ColumnWidth = 18.14
$ws.Range("C1").EntireColumn.ColumnWidth = 15.63
